$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '43.828.67'
Set-TextValue 'E2' '  +4.46%  '
Set-TextValue 'D3' '2.274.93'
Set-TextValue 'E3' '  +1.99%  '
Set-TextValue 'E4' '  +0.01%  '
Set-TextValue 'D5' '231.30'
Set-TextValue 'E5' '  -0.14%  '
Set-TextValue 'D6' '0.628'
Set-TextValue 'E6' '  +0.57%  '
Set-TextValue 'D7' '61.49'
Set-TextValue 'E7' '  +0.77%  '
Set-TextValue 'E8' '  +0.02%  '
Set-TextValue 'D9' '0.423'
Set-TextValue 'E9' '  +5.54%  '
Set-TextValue 'D10' '0.0944'
Set-TextValue 'E10' '  +6.09%  '
Set-TextValue 'D11' '57.80'
Set-TextValue 'E11' '  -1.97%  '
Set-TextValue 'E12' '  +0.77%  '
Set-TextValue 'D13' '2.613.41'
Set-TextValue 'E13' '  +2.10%  '
Set-TextValue 'D14' '15.79'
Set-TextValue 'E14' '  +1.00%  '
Set-TextValue 'D15' '23.72'
Set-TextValue 'E15' '  +9.13%  '
Set-TextValue 'D16' '5.81'
Set-TextValue 'E16' '  +4.17%  '
Set-TextValue 'D17' '0.812'
Set-TextValue 'E17' '  +1.76%  '
Set-TextValue 'D18' '2.276.34'
Set-TextValue 'E18' '  +1.24%  '
Set-TextValue 'D19' '43.751.96'
Set-TextValue 'E19' '  +4.67%  '
Set-TextValue 'E20' '  +5.50%  '
Set-TextValue 'D21' '73.14'
Set-TextValue 'E21' '  +0.90%  '
Set-TextValue 'D22' '6.24'
Set-TextValue 'E22' '  +3.44%  '
Set-TextValue 'D23' '251.76'
Set-TextValue 'E23' '  +0.70%  '
Set-TextValue 'E24' '  +0.06%  '
Set-TextValue 'E25' '  +7.49%  '
Set-TextValue 'E26' '  +2.36%  '
Set-TextValue 'D27' '9.87'
Set-TextValue 'E27' '  +1.87%  '
Set-TextValue 'D28' '171.07'
Set-TextValue 'E28' '  +2.11%  '
Set-TextValue 'E29' '  -1.31%  '
Set-TextValue 'D30' '20.60'
Set-TextValue 'E30' '  +3.30%  '
Set-TextValue 'D31' '1.46'
Set-TextValue 'E31' '  +4.30%  '
Set-TextValue 'E32' '  +1.08%  '
Set-TextValue 'E33' '  +0.05%  '
Set-TextValue 'D34' '4.80'
Set-TextValue 'E34' '  +3.90%  '
Set-TextValue 'D36' '0.0662'
Set-TextValue 'E36' '  +5.19%  '
Set-TextValue 'D37' '6.49'
Set-TextValue 'E37' '  -2.18%  '
Set-TextValue 'D38' '2.41'
Set-TextValue 'E38' '  +2.15%  '
Set-TextValue 'D39' '3.61'
Set-TextValue 'E39' '  -1.91%  '
Set-TextValue 'E40' '  +4.66%  '
Set-TextValue 'D41' '1.00'
Set-TextValue 'E41' '  -0.09%  '
Set-TextValue 'D42' '8.77'
Set-TextValue 'E42' '  +2.47%  '
Set-TextValue 'D43' '0.000225'
Set-TextValue 'E43' '  -11.95%  '
Set-TextValue 'B44' 'Cronos'
Set-TextValue 'C44' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D44' '0.0986'
Set-TextValue 'E44' '  +0.64%  '
Set-TextValue 'B45' 'FTXToken'
Set-TextValue 'C45' 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue 'D45' '4.55'
Set-TextValue 'E45' '  -5.21%  '
Set-TextValue 'E46' '  +0.66%  '
Set-TextValue 'D47' '98.02'
Set-TextValue 'E47' '  -0.82%  '
Set-TextValue 'D48' '1.470.39'
Set-TextValue 'E48' '  -0.04%  '
Set-TextValue 'D49' '16.64'
Set-TextValue 'E49' '  +1.00%  '
Set-TextValue 'E50' '  +0.86%  '
Set-TextValue 'B51' 'NEARProtocol'
Set-TextValue 'C51' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D51' '2.27'
Set-TextValue 'E51' '  +8.76%  '
